$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: wrap a bare run-level WordprocessingML fragment (<w:r>.../<w:r>,
# <w:bookmarkStart/>, <w:proofErr/>, etc.) in the full pkg:package envelope
# that InsertXML expects, then apply it to a Range, replacing that range's
# contents with the supplied runs.
# ---------------------------------------------------------------------------
function Apply-RunXml($range, [string]$innerRuns) {
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p>' + $innerRuns + '</w:p></w:body></w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'
    $range.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# 1) Drop the "_GoBack" bookmark from its old spot (start of the "Glossary
#    and formulas..." paragraph) - it is being relocated into the title.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2) Title paragraph: split "Scoring Guidance Principles Based on Attribute
#    Scores" into four runs ("Priority S" / "coring " / "of " /
#    "Guidance Principles Based on Attribute Scores") and drop the
#    relocated "_GoBack" bookmark right before the last run.
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs(1).Range
$titleRange = $d.Range($titlePara.Start, $titlePara.End - 1)

$titleRuns =
    '<w:r><w:rPr><w:sz w:val="32"/><w:u w:val="single"/></w:rPr><w:t>Priority S</w:t></w:r>' +
    '<w:r><w:rPr><w:sz w:val="32"/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve">coring </w:t></w:r>' +
    '<w:r><w:rPr><w:sz w:val="32"/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve">of </w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
    '<w:r><w:rPr><w:sz w:val="32"/><w:u w:val="single"/></w:rPr><w:t>Guidance Principles Based on Attribute Scores</w:t></w:r>'

Apply-RunXml $titleRange $titleRuns

# ---------------------------------------------------------------------------
# 3) "Glossary and formulas..." paragraph: append a new sentence describing
#    the value range, right after the closing ")".
# ---------------------------------------------------------------------------
$glossaryPara = $d.Paragraphs(3).Range
$glossaryRange = $d.Range($glossaryPara.Start, $glossaryPara.End - 1)

$glossaryRuns =
    '<w:r w:rsidRPr="00F14B8F"><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Glossary and formulas for the guiding principles. Each guiding principle </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r w:rsidRPr="00F14B8F"><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>is derived</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r w:rsidRPr="00F14B8F"><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> from the scores from the limiting attributes described in Table 1. </w:t></w:r>' +
    '<w:r w:rsidRPr="00F14B8F"><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/><w:b/><w:bCs/><w:i/><w:iCs/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Data </w:t></w:r>' +
    '<w:r w:rsidRPr="00F14B8F"><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">indicates all data-limitations (i.e., # Types, Precision, Bias, Species ID, Spatial, Temporal). </w:t></w:r>' +
    '<w:r w:rsidRPr="00F14B8F"><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/><w:b/><w:bCs/><w:i/><w:iCs/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Resource</w:t></w:r>' +
    '<w:r w:rsidRPr="00F14B8F"><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> indicates all resource-limitations (i.e., Time, Funding, Capacity, </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r w:rsidRPr="00F14B8F"><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Analysts</w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r w:rsidRPr="00F14B8F"><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>:Stocks</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/>' +
    '<w:r w:rsidRPr="00F14B8F"><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>)</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">. Values range from </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>0</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> to 3, with 0 being the lowest priority score and 3 being the highest priority score.</w:t></w:r>'

Apply-RunXml $glossaryRange $glossaryRuns

Write-Output "Done."
